$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style from an untouched default-style cell, used to restore
# plain (no explicit format) styling on numeric-looking text cells below,
# since assigning a plain numeric-looking string to .Value would otherwise
# be auto-converted to a number by Excel.
$defaultStyle = $ws.Range("D4").Style

$ws.Range("D2").Value = '27.662.20'
$ws.Range("D3").Value = '1.583.28'
$ws.Range("E3").Value = '  -3.27%  '
$ws.Range("E4").Value = '  +0.36%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.12'
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = '  -2.69%  '
$ws.Range("E6").Value = '  -2.90%  '
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.23'
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = '  -5.06%  '
$ws.Range("E9").Value = '  -1.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0591'
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = '  -3.27%  '
$ws.Range("E11").Value = '  -1.97%  '
$ws.Range("D12").Value = '1.809.86'
$ws.Range("E12").Value = '  -3.17%  '
$ws.Range("D13").Value = '1.583.49'
$ws.Range("E13").Value = '  -3.23%  '
$ws.Range("E14").Value = '  -4.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.529'
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = '  -6.22%  '
$ws.Range("D16").Value = '27.646.21'
$ws.Range("E16").Value = '  -1.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.26'
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = '  -3.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '219.96'
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = '  -4.16%  '
$ws.Range("D19").Value = '0.0₃0691'
$ws.Range("E19").Value = '  -3.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.29'
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = '  -6.59%  '
$ws.Range("E21").Value = '  +0.35%  '
$ws.Range("E22").Value = '  -5.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.47'
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = '  -6.96%  '
$ws.Range("E24").Value = '  -5.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.98'
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = '  -1.47%  '
$ws.Range("E26").Value = '  -2.86%  '
$ws.Range("E27").Value = '  +0.26%  '
$ws.Range("E28").Value = '  -2.85%  '
$ws.Range("E29").Value = '  -4.12%  '
$ws.Range("E30").Value = '  -2.70%  '
$ws.Range("E31").Value = '  -3.74%  '
$ws.Range("E32").Value = '  -6.00%  '
$ws.Range("D33").Value = '1.384.48'
$ws.Range("E33").Value = '  -1.36%  '
$ws.Range("E34").Value = '  -5.80%  '
$ws.Range("E35").Value = '  -5.54%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.964'
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = '  -5.34%  '
$ws.Range("E37").Value = '  -0.88%  '
$ws.Range("E38").Value = '  -3.34%  '
$ws.Range("E39").Value = '  -3.59%  '
$ws.Range("E40").Value = '  -4.16%  '
$ws.Range("E41").Value = '  +0.37%  '
$ws.Range("E42").Value = '  -2.98%  '
$ws.Range("E43").Value = '  -4.43%  '
$ws.Range("E44").Value = '  +1.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.48'
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = '  -4.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.22'
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = '  -4.55%  '
$ws.Range("D47").Value = '1.721.08'
$ws.Range("E47").Value = '  -3.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.83'
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = '  -1.01%  '
$ws.Range("E49").Value = '  -2.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0971'
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = '  -5.17%  '
$ws.Range("E51").Value = '  -1.12%  '
